# Update the "取得日時" (acquisition timestamp) column for all data rows
# on the "ランサーズ" sheet to reflect the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-25 06:38:42"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
